$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.886.63"
$ws.Range("E2").Value = "  +3.34%  "
$ws.Range("D3").Value = "3.572.95"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'582.30"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").Value = "'186.68"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  +2.49%  "
$ws.Range("D8").Value = "3.560.98"
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "'0.224"
$ws.Range("E10").Value = "  +23.12%  "
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").Value = "'54.68"
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("D13").Value = "'0.0000320"
$ws.Range("E13").Value = "  +6.93%  "
$ws.Range("D14").Value = "'9.48"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").Value = "4.140.70"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "70.891.22"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").Value = "'19.23"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "'12.82"
$ws.Range("E18").Value = "  +5.25%  "
$ws.Range("D19").Value = "3.564.60"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("D20").Value = "'574.79"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "'17.62"
$ws.Range("E23").Value = "  -7.63%  "
$ws.Range("E24").Value = "  +4.64%  "
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("D26").Value = "'94.02"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "'11.24"
$ws.Range("E27").Value = "  +4.94%  "
$ws.Range("D28").Value = "'2.96"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("D29").Value = "'9.23"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").Value = "'32.45"
$ws.Range("E30").Value = "  +3.43%  "
$ws.Range("D31").Value = "'7.21"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").Value = "'63.19"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("D35").Value = "'3.40"
$ws.Range("E35").Value = "  +14.87%  "
$ws.Range("D36").Value = "'3.60"
$ws.Range("E36").Value = "  +17.23%  "
$ws.Range("D37").Value = "'548.69"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("E38").Value = "  +5.55%  "
$ws.Range("D39").Value = "'38.19"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "0.0₃0803"
$ws.Range("E40").Value = "  +5.56%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "3.580.94"
$ws.Range("E42").Value = "  +10.38%  "
$ws.Range("E43").Value = "  +5.11%  "
$ws.Range("D44").Value = "'3.44"
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("D45").Value = "'0.0469"
$ws.Range("E45").Value = "  +7.78%  "
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "'2.92"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  +4.83%  "
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("E50").Value = "  +14.58%  "
$ws.Range("E51").Value = "  -0.04%  "
